{"js": "// Update the TINKER license header:\n//   \"Version 6.0   October 2011\"  -> \"Version 6.2   February 2013\"\n//   \"Copyright \u00a9 1990-2011\"       -> \"Copyright \u00a9 1990-2013\"\n// and relocate the \"_GoBack\" bookmark so it sits at the start of the\n// Copyright paragraph (matching the author's edit, which moved the\n// paragraph break to just before the bookmark).\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) \"Version 6.0   October 2011\" -> \"Version 6.2   February 2013\"\nlet results = body.search(\"Version 6.0   October 2011\", { matchCase: true });\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Version 6.2   February 2013\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"Copyright \u00a9 1990-2011\" -> \"Copyright \u00a9 1990-2013\"\nresults = body.search(\"Copyright \u00a9 1990-2011\", { matchCase: true });\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Copyright \u00a9 1990-2013\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Move the \"_GoBack\" bookmark from the version line to the start of the\n//    Copyright line (the diff shows the paragraph break now lands right\n//    before the bookmark).\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nresults = body.search(\"Copyright \u00a9 1990-2013\", { matchCase: true });\nawait context.sync();\nif (results.items.length > 0) {\n  const startRange = results.items[0].getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Update the TINKER license header:\n#   \"Version 6.0   October 2011\"  -> \"Version 6.2   February 2013\"\n#   \"Copyright \u00a9 1990-2011\"       -> \"Copyright \u00a9 1990-2013\"\n# and relocate the \"_GoBack\" bookmark so it sits at the start of the\n# Copyright paragraph (matching the author's edit, which moved the\n# paragraph break to just before the bookmark).\n\n$d = $word.ActiveDocument\n\n# Remove the existing \"_GoBack\" bookmark (it sits inside the text we are\n# about to replace, so Word would drop it anyway) - we re-add it below at\n# its new location.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 1) \"Version 6.0   October 2011\" -> \"Version 6.2   February 2013\"\n$r1 = $d.Content\n$f1 = $r1.Find\n$f1.ClearFormatting()\n$f1.Text = \"Version 6.0   October 2011\"\n$found1 = $f1.Execute()\nif ($found1) {\n    $r1.Text = \"Version 6.2   February 2013\"\n}\n\n# 2) \"Copyright \u00a9 1990-2011\" -> \"Copyright \u00a9 1990-2013\"\n$r2 = $d.Content\n$f2 = $r2.Find\n$f2.ClearFormatting()\n$f2.Text = \"Copyright \u00a9 1990-2011\"\n$found2 = $f2.Execute()\nif ($found2) {\n    $r2.Text = \"Copyright \u00a9 1990-2013\"\n}\n\n# 3) Re-insert \"_GoBack\" at the start of the (now updated) Copyright line.\n$r3 = $d.Content\n$f3 = $r3.Find\n$f3.ClearFormatting()\n$f3.Text = \"Copyright \u00a9 1990-2013\"\n$found3 = $f3.Execute()\nif ($found3) {\n    $bmRange = $r3.Duplicate\n    $bmRange.Collapse(1)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
